$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1199.875
$ws.Range("I40").Value = 1180
$ws.Range("J40").Value = 1233
$ws.Range("K40").Value = 1180
$ws.Range("L40").Value = 1233
$ws.Range("M40").Value = -1005
$ws.Range("N40").Value = -1583
$ws.Range("H76").Value = 2868.5715
$ws.Range("I76").Value = 2826.75
$ws.Range("J76").Value = 3705
$ws.Range("K76").Value = 2826.75
$ws.Range("L76").Value = 3705
$ws.Range("M76").Value = -2511.75
$ws.Range("N76").Value = -4335
$ws.Range("H79").Value = 2868.5715
$ws.Range("I79").Value = 2826.75
$ws.Range("J79").Value = 3705
$ws.Range("K79").Value = 2826.75
$ws.Range("L79").Value = 3705
$ws.Range("M79").Value = -1734.75
$ws.Range("N79").Value = -5889
$ws.Range("H132").Value = 2567018.8
$ws.Range("I132").Value = 2900117.5
$ws.Range("K132").Value = 8700352.5
$ws.Range("M132").Value = -8697822.5
$ws.Range("H137").Value = 20849.537
$ws.Range("I137").Value = 2443.205
$ws.Range("K137").Value = 7329.615
$ws.Range("M137").Value = -4779.615
$ws.Range("H138").Value = 3926.5386
$ws.Range("J138").Value = 5106.8237
$ws.Range("L138").Value = 15320.4711
$ws.Range("N138").Value = -25600.4711

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 26084.334
$ws.Range("I23").Value = 70006
$ws.Range("J23").Value = 17300
$ws.Range("K23").Value = 70006
$ws.Range("L23").Value = 17300
$ws.Range("M23").Value = -69747
$ws.Range("N23").Value = -17818
$ws.Range("H32").Value = 1631.4
$ws.Range("I32").Value = 1418.4667
$ws.Range("J32").Value = 3547.8
$ws.Range("K32").Value = 1418.4667
$ws.Range("L32").Value = 3547.8
$ws.Range("M32").Value = -1131.4667
$ws.Range("N32").Value = -4121.8
$ws.Range("H61").Value = 2334.16
$ws.Range("I61").Value = 931.5
$ws.Range("J61").Value = 4827.778
$ws.Range("K61").Value = 931.5
$ws.Range("L61").Value = 4827.778
$ws.Range("M61").Value = -719.5
$ws.Range("N61").Value = -5251.778
$ws.Range("H63").Value = 3243.75
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 3243.75
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H122").Value = 3338.1765
$ws.Range("I122").Value = 2520.8333
$ws.Range("J122").Value = 5299.8
$ws.Range("K122").Value = 7562.499899999999
$ws.Range("L122").Value = 15899.4
$ws.Range("M122").Value = -5112.499899999999
$ws.Range("N122").Value = -20799.4
$ws.Range("H136").Value = 2334.16
$ws.Range("I136").Value = 931.5
$ws.Range("J136").Value = 4827.778
$ws.Range("K136").Value = 2794.5
$ws.Range("L136").Value = 14483.334
$ws.Range("M136").Value = -244.5
$ws.Range("N136").Value = -19583.334

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 5000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""
$ws.Range("H82").Value = 15657.5
$ws.Range("J82").Value = 38848.6
$ws.Range("L82").Value = 38848.6
$ws.Range("N82").Value = -39614.6
$ws.Range("H85").Value = 15657.5
$ws.Range("J85").Value = 38848.6
$ws.Range("L85").Value = 38848.6
$ws.Range("N85").Value = -41500.6

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 5666.6665
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 6000
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = -4860
$ws.Range("N8").Value = -6280
$ws.Range("H31").Value = 181161.9
$ws.Range("I31").Value = 1642.4193
$ws.Range("K31").Value = 1642.4193
$ws.Range("M31").Value = -1347.4193
$ws.Range("H34").Value = 181161.9
$ws.Range("I34").Value = 1642.4193
$ws.Range("K34").Value = 1642.4193
$ws.Range("M34").Value = -1440.4193
$ws.Range("H41").Value = 7884.6665
$ws.Range("I41").Value = 2500
$ws.Range("J41").Value = 10577
$ws.Range("K41").Value = 2500
$ws.Range("L41").Value = 10577
$ws.Range("M41").Value = -2072
$ws.Range("N41").Value = -11433
$ws.Range("H50").Value = 18092
$ws.Range("J50").Value = 18092
$ws.Range("L50").Value = 18092
$ws.Range("N50").Value = -19342
$ws.Range("H51").Value = 9166.333000000001
$ws.Range("J51").Value = 9166.333000000001
$ws.Range("L51").Value = 9166.333000000001
$ws.Range("N51").Value = -10638.333
$ws.Range("H60").Value = 16400.75
$ws.Range("J60").Value = 16534.334
$ws.Range("L60").Value = 16534.334
$ws.Range("N60").Value = -17556.334
$ws.Range("H61").Value = 9166.333000000001
$ws.Range("J61").Value = 9166.333000000001
$ws.Range("L61").Value = 9166.333000000001
$ws.Range("N61").Value = -9862.333000000001
$ws.Range("H68").Value = 35970
$ws.Range("J68").Value = 35970
$ws.Range("L68").Value = 35970
$ws.Range("N68").Value = -37468
$ws.Range("H71").Value = 35970
$ws.Range("J71").Value = 35970
$ws.Range("L71").Value = 107910
$ws.Range("N71").Value = -115398
$ws.Range("H74").Value = 16132
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 16132
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -17880
$ws.Range("H77").Value = 16132
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 16132
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -57132
$ws.Range("H132").Value = 1900.9423
$ws.Range("I132").Value = 1396.619
$ws.Range("J132").Value = 4019.1
$ws.Range("K132").Value = 4189.857
$ws.Range("L132").Value = 12057.3
$ws.Range("M132").Value = -1659.857
$ws.Range("N132").Value = -17117.3

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2420
$ws.Range("I57").Value = 840
$ws.Range("K57").Value = 2520
$ws.Range("M57").Value = -1961
$ws.Range("H69").Value = 86443.234
$ws.Range("I69").Value = 620.6667
$ws.Range("K69").Value = 1862.0001
$ws.Range("M69").Value = -1051.0001
$ws.Range("H72").Value = 86443.234
$ws.Range("I72").Value = 620.6667
$ws.Range("K72").Value = 5586.0003
$ws.Range("M72").Value = -1530.0003

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2349.4082
$ws.Range("I132").Value = 1516
$ws.Range("J132").Value = 3291.5217
$ws.Range("K132").Value = 4548
$ws.Range("L132").Value = 9874.5651
$ws.Range("M132").Value = -2018
$ws.Range("N132").Value = -14934.5651
$ws.Range("H136").Value = 1528.8334
$ws.Range("I136").Value = 949.6896400000001
$ws.Range("J136").Value = 3928.1428
$ws.Range("K136").Value = 2849.06892
$ws.Range("L136").Value = 11784.4284
$ws.Range("M136").Value = -299.0689200000002
$ws.Range("N136").Value = -16884.4284
